# Apply the updated crypto price/volume figures to Sheet1 (rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.863.17"
$ws.Range("E2").Value = "  -1.14%  "

$ws.Range("D3").Value = "1.563.57"
$ws.Range("E3").Value = "  +0.25%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.93"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.32%  "

$ws.Range("E6").Value = "  -1.53%  "

$ws.Range("E7").Value = "  -0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.82"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.89%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.247"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.54%  "

$ws.Range("E10").Value = "  -1.39%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0865"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.39%  "

$ws.Range("D12").Value = "1.784.56"
$ws.Range("E12").Value = "  +0.23%  "

$ws.Range("D13").Value = "1.578.24"
$ws.Range("E13").Value = "  +1.22%  "

$ws.Range("E14").Value = "  -1.33%  "

$ws.Range("E15").Value = "  -0.50%  "

$ws.Range("D16").Value = "26.863.53"
$ws.Range("E16").Value = "  -1.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.28"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "215.20"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.59%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.36"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.77%  "

$ws.Range("D20").Value = "0.0₃0681"
$ws.Range("E20").Value = "  -0.92%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.13"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.30%  "

$ws.Range("E23").Value = "  -1.83%  "

$ws.Range("E24").Value = "  +1.35%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.76"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.09%  "

$ws.Range("E26").Value = "  +1.56%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.92"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.20%  "

$ws.Range("E28").Value = "  -0.15%  "

$ws.Range("E29").Value = "  -1.33%  "

$ws.Range("E30").Value = "  +0.98%  "

$ws.Range("E31").Value = "  -3.78%  "

$ws.Range("E32").Value = "  -0.14%  "

$ws.Range("D33").Value = "1.400.93"
$ws.Range("E33").Value = "  +1.23%  "

$ws.Range("E34").Value = "  -0.45%  "

$ws.Range("E35").Value = "  -1.19%  "

$ws.Range("E36").Value = "  -0.48%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.910"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.01%  "

$ws.Range("E38").Value = "  -1.05%  "

$ws.Range("E39").Value = "  +2.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.811"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.25%  "

$ws.Range("E41").Value = "  -0.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.996"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.14%  "

$ws.Range("E43").Value = "  +6.35%  "

$ws.Range("E44").Value = "  -0.51%  "

$ws.Range("E45").Value = "  +1.05%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.41"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.29%  "

$ws.Range("D47").Value = "1.698.84"
$ws.Range("E47").Value = "  +0.32%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "86.55"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0502"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.87%  "

$ws.Range("D50").Value = "0.0₇0968"
$ws.Range("E50").Value = "  -1.42%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0949"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.80%  "
